# Ajout curseurs regexpr parties variables pour 2019
#
# Duplicates the existing 2018 "curseur" block (rows 68:75 of Feuil1) as a
# new 2019 block (rows 76:83), reusing the same table/zone/variable labels
# but with the year bumped to 2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$data = @(
    @("rum", 2019, "zac",  "[A-Z]{4}[0-9]{3}", 7),
    @("rum", 2019, "zd",   ".{1,8}",           8),
    @("rum", 2019, "zdad", ".{1,8}",           8),
    @("rum", 2019, "zal",  ".{1,29}",          29),
    @("rsa", 2019, "zac",  "[A-Z]{4}[0-9]{3}", 7),
    @("rsa", 2019, "zd",   ".{1,6}",           6),
    @("rsa", 2019, "zum",  ".{1,60}",          60),
    @("rsa", 2019, "zal",  ".{24}",            24)
)

$startRow = 76
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Match the author's final selection (B76:B83) shown in the saved sheetView.
$ws.Range("B76:B83").Select() | Out-Null

# Leftover-from-session AutoFilter defined name the saved workbook carries
# (sheet-scoped, hidden) - mirrors what Excel stamps when a filter has been
# toggled on the sheet.
try {
    $n = $ws.Names.Add("_xlnm._FilterDatabase", "=Feuil1!`$A`$1:`$E`$75")
    $n.Visible = $false
} catch {
}
